$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Description column (B) text for rows 2-23
$ws.Cells.Item(2, 2).Value = 'Ham, provolone, lettuce, tomatoes, onions, and mayo'
$ws.Cells.Item(3, 2).Value = 'Turkey, provolone, lettuce, tomatoes, onions, and mayo'
$ws.Cells.Item(4, 2).Value = 'Roast beef, 2x provolone, lettuce, tomatoes, onions, and Roasted Garlic Aioli'
$ws.Cells.Item(5, 2).Value = 'Pepperoni, salami, turkey, ham and roast beef, 2x provolone, lettuce, tomatoes, onions, mayo, and MVP Vinaigrette'
$ws.Cells.Item(6, 2).Value = 'Steak, 2x American Cheese, Green Peppers, Red Onions, Toasted, and Mayo'
$ws.Cells.Item(7, 2).Value = 'Steak, 2x American Cheese, Green Peppers, Red Onions, Toasted, and Sweet Onion Teriyaki sauce'
$ws.Cells.Item(8, 2).Value = 'Steak, 2x pepperjack cheese, green peppers, red onions, baja chiptole sauce, Toasted'
$ws.Cells.Item(9, 2).Value = 'Steak, bacon, Monterey cheddar, green peppers and red onions piled high, Artisan Italian bread, creamy Peppercorn Ranch'
$ws.Cells.Item(10, 2).Value = 'Genoa salami, spicy pepperoni, jalapeno peppers, lettuce, tomato, red onions, provolone cheese, and MVP Parmesan Vinaigrette'
$ws.Cells.Item(11, 2).Value = 'Genoa salami, spicy pepperoni, savory Black Forest ham, provolone cheese, crisp lettuce, tomatoes, red onions'
$ws.Cells.Item(12, 2).Value = 'Black Forest ham, Genoa salami, pepperoni, capicola on Italian bread with provolone cheese, lettuce, tomatoes, red onions, and banana peppers'
$ws.Cells.Item(13, 2).Value = 'Thin-sliced Black Forest ham, capicola, and BelGioioso® Fresh Mozzarella on Italian bread; with spinach, tomatoes, red onions, and banana peppers'
$ws.Cells.Item(14, 2).Value = 'Juicy meatballs drenched in marinara sauce with slices of pepperoni and BelGioioso® Fresh Mozzarella'
$ws.Cells.Item(15, 2).Value = 'Juicy rotisserie-style chicken, crispy bacon, provolone, lettuce, tomatoes, red onions '
$ws.Cells.Item(16, 2).Value = 'Multigrain bread, grilled chicken strips marinated in Sweet Onion Teriyaki sauce, American cheese, lettuce, tomatoes, and red onions.'
$ws.Cells.Item(17, 2).Value = 'Tender rotisserie-style chicken, melted Monterey cheddar cheese, hickory-smoked bacon, lettuce, tomatoes, red onions,'
$ws.Cells.Item(18, 2).Value = 'Tender hand-pulled rotisserie-style chicken, Monterey Cheddar, green peppers, and red onions'
$ws.Cells.Item(19, 2).Value = 'Juicy rotisserie-style chicken, smashed avocado, double Pepper Jack cheese, lettuce, tomatoes, and red onions'
$ws.Cells.Item(20, 2).Value = 'Oven-roasted turkey, Black Forest ham, crisp bacon, and American cheese with lettuce, tomatoes, red onions, and mayo'
$ws.Cells.Item(21, 2).Value = 'Oven-roasted turkey, Black Forest ham, crisp bacon, and American cheese with lettuce, tomatoes, red onions, and mayo'
$ws.Cells.Item(22, 2).Value = 'Oven-roasted turkey, Black Forest ham, Roast Beef, provolone cheese on Multigrain bread topped with lettuce, tomatoes, red onions, and mayo'
$ws.Cells.Item(23, 2).Value = 'Oven-roasted turkey, crisp bacon, BelGioioso® Fresh Mozzarella, and smashed Hass avocado. Topped with spinach, tomatoes, red onions, and mayo'

# Update row heights
$ws.Rows(2).RowHeight = 31
$ws.Rows(3).RowHeight = 31
$ws.Rows(4).RowHeight = 46.5
$ws.Rows(5).RowHeight = 62
$ws.Rows(6).RowHeight = 46.5
$ws.Rows(7).RowHeight = 46.5
$ws.Rows(8).RowHeight = 46.5
$ws.Rows(9).RowHeight = 62
$ws.Rows(10).RowHeight = 77.5
$ws.Rows(11).RowHeight = 62
$ws.Rows(12).RowHeight = 77.5
$ws.Rows(13).RowHeight = 77.5
$ws.Rows(14).RowHeight = 62
$ws.Rows(15).RowHeight = 46.5
$ws.Rows(16).RowHeight = 77.5
$ws.Rows(17).RowHeight = 31
$ws.Rows(18).RowHeight = 46.5
$ws.Rows(19).RowHeight = 62
$ws.Rows(20).RowHeight = 62
$ws.Rows(21).RowHeight = 62
$ws.Rows(22).RowHeight = 77.5
$ws.Rows(23).RowHeight = 77.5

# Update sheet view: set new selection (also resets scrolled top-left cell)
$null = $ws.Range("F5").Select()

Write-Output "done"
